# Apply the "Add files via upload" edit:
#  - rename the sheet from "Sheet1" to "Main"
#  - highlight the title cell (A1) and the "Unaccounted" total (B24) with a
#    light gold fill (Excel theme color "Gold, Accent 4, Lighter 80%",
#    i.e. theme index 7 with tint ~0.8 -> RGB #FFF2CC)
#  - leave the cursor on A1 (matches the saved sheetView in the source file)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet tab.
$ws.Name = "Main"

# 2. Fill B24 ("Unaccounted" value) first, then A1 (report title) second, so
#    the new cellXfs entries are appended in the same order seen in the
#    target workbook (B24's style before A1's style).
$ws.Range("B24").Interior.Color = 13431551   # RGB(255,242,204) == #FFF2CC
$ws.Range("A1").Interior.Color = 13431551    # RGB(255,242,204) == #FFF2CC

# 3. Leave the selection on A1 (title cell), matching the saved view state.
$ws.Range("A1").Select() | Out-Null
